$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "60.917.65"
$ws.Range("E2").Value2 = "  -0.13%  "

$ws.Range("D3").Value2 = "2.911.39"
$ws.Range("E3").Value2 = "  -0.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value2 = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.05"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value2 = "  +0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.82"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value2 = "  -0.47%  "

$ws.Range("E7").Value2 = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.505"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value2 = "  -0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.89"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value2 = "  +1.22%  "

$ws.Range("E10").Value2 = "  -2.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.439"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value2 = "  -2.12%  "

$ws.Range("E12").Value2 = "  -0.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.39"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value2 = "  -0.63%  "

$ws.Range("E14").Value2 = "  -0.05%  "

$ws.Range("D15").Value2 = "3.393.59"
$ws.Range("E15").Value2 = "  -0.33%  "

$ws.Range("D16").Value2 = "60.801.34"
$ws.Range("E16").Value2 = "  -0.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.67"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value2 = "  -1.30%  "

$ws.Range("D18").Value2 = "2.907.14"
$ws.Range("E18").Value2 = "  -0.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "432.16"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value2 = "  +0.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.35"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value2 = "  -1.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.677"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value2 = "  -0.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.09"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value2 = "  -0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.43"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value2 = "  +0.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.85"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value2 = "  +0.77%  "

$ws.Range("E25").Value2 = "  -2.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.79"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value2 = "  -2.04%  "

$ws.Range("E27").Value2 = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value2 = "  +4.66%  "

$ws.Range("E29").Value2 = "  -1.10%  "

$ws.Range("E30").Value2 = "  -4.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.50"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value2 = "  -0.56%  "

$ws.Range("E32").Value2 = "  +1.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value2 = "  -0.16%  "

$ws.Range("D34").Value2 = "0.0₃0860"
$ws.Range("E34").Value2 = "  -1.35%  "

$ws.Range("E35").Value2 = "  -0.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.61"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value2 = "  -0.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.99"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value2 = "  -1.16%  "

$ws.Range("E38").Value2 = "  -1.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.121"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value2 = "  -4.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.54"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value2 = "  -1.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "41.25"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value2 = "  +0.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.282"
$ws.Range("D42").NumberFormat = "General"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "375.46"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value2 = "  -1.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0344"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value2 = "  -2.52%  "

$ws.Range("D45").Value2 = "2.697.33"
$ws.Range("E45").Value2 = "  +0.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.56"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value2 = "  +0.67%  "

$ws.Range("E47").Value2 = "  -0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.72"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value2 = "  -3.73%  "

$ws.Range("E49").Value2 = "  -0.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.00"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value2 = "  -3.33%  "

$ws.Range("E51").Value2 = "  -1.19%  "
